$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 8507.333000000001
$ws.Range("I53").Value = 75012.5
$ws.Range("J53").Value = 194.1875
$ws.Range("K53").Value = 75012.5
$ws.Range("L53").Value = 194.1875
$ws.Range("M53").Value = -74375.5
$ws.Range("N53").Value = -1468.1875

$ws.Range("H70").Value = 7722.077
$ws.Range("I70").Value = 1075
$ws.Range("J70").Value = 8930.637000000001
$ws.Range("K70").Value = 3225
$ws.Range("L70").Value = 26791.911
$ws.Range("M70").Value = -2955
$ws.Range("N70").Value = -27331.911

$ws.Range("H73").Value = 7722.077
$ws.Range("I73").Value = 1075
$ws.Range("J73").Value = 8930.637000000001
$ws.Range("K73").Value = 3225
$ws.Range("L73").Value = 26791.911
$ws.Range("M73").Value = -2289
$ws.Range("N73").Value = -28663.911

$ws.Range("H106").Value = 34484536
$ws.Range("I106").Value = 37038804
$ws.Range("K106").Value = 37038804
$ws.Range("M106").Value = -37038173

$ws.Range("H115").Value = 7722.85
$ws.Range("I115").Value = 696.7143
$ws.Range("K115").Value = 2090.1429
$ws.Range("M115").Value = -523.1428999999998

$ws.Range("H125").Value = 1872041.4
$ws.Range("I125").Value = 1476
$ws.Range("J125").Value = 3041144.8
$ws.Range("K125").Value = 13284
$ws.Range("L125").Value = 27370303.2
$ws.Range("M125").Value = -10824
$ws.Range("N125").Value = -27375223.2

$ws.Range("H137").Value = 15256862
$ws.Range("I137").Value = 1012.65515
$ws.Range("J137").Value = 52125170
$ws.Range("K137").Value = 3037.96545
$ws.Range("L137").Value = 156375510
$ws.Range("M137").Value = -487.9654500000001
$ws.Range("N137").Value = -156380610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 13908.909
$ws.Range("J55").Value = 13908.909
$ws.Range("L55").Value = 13908.909
$ws.Range("N55").Value = -14538.909

$ws.Range("H74").Value = 76669000
$ws.Range("I74").Value = 100001336
$ws.Range("J74").Value = 53336670
$ws.Range("K74").Value = 100001336
$ws.Range("L74").Value = 53336670
$ws.Range("M74").Value = -100000462
$ws.Range("N74").Value = -53338418

$ws.Range("H77").Value = 76669000
$ws.Range("I77").Value = 100001336
$ws.Range("J77").Value = 53336670
$ws.Range("K77").Value = 500006680
$ws.Range("L77").Value = 266683350
$ws.Range("M77").Value = -500002312
$ws.Range("N77").Value = -266692086

$ws.Range("H80").Value = 24286
$ws.Range("J80").Value = 24286
$ws.Range("L80").Value = 24286
$ws.Range("N80").Value = -26282

$ws.Range("H83").Value = 24286
$ws.Range("J83").Value = 24286
$ws.Range("L83").Value = 72858
$ws.Range("N83").Value = -82842

$ws.Range("H102").Value = 1394.0625
$ws.Range("I102").Value = 1394.0625
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1394.0625
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 227.9375
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18465.312
$ws.Range("I82").Value = 1918.6
$ws.Range("J82").Value = 25986.545
$ws.Range("K82").Value = 1918.6
$ws.Range("L82").Value = 25986.545
$ws.Range("M82").Value = -1535.6
$ws.Range("N82").Value = -26752.545

$ws.Range("H85").Value = 18465.312
$ws.Range("I85").Value = 1918.6
$ws.Range("J85").Value = 25986.545
$ws.Range("K85").Value = 1918.6
$ws.Range("L85").Value = 25986.545
$ws.Range("M85").Value = -592.5999999999999
$ws.Range("N85").Value = -28638.545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1457183
$ws.Range("I31").Value = 1735.8462
$ws.Range("J31").Value = 3683161
$ws.Range("K31").Value = 1735.8462
$ws.Range("L31").Value = 3683161
$ws.Range("M31").Value = -1440.8462
$ws.Range("N31").Value = -3683751

$ws.Range("H34").Value = 1457183
$ws.Range("I34").Value = 1735.8462
$ws.Range("J34").Value = 3683161
$ws.Range("K34").Value = 1735.8462
$ws.Range("L34").Value = 3683161
$ws.Range("M34").Value = -1533.8462
$ws.Range("N34").Value = -3683565

$ws.Range("H51").Value = 9455.777
$ws.Range("J51").Value = 8887.75
$ws.Range("L51").Value = 8887.75
$ws.Range("N51").Value = -10359.75

$ws.Range("H59").Value = 15571.857
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 15571.857
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 15571.857
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -17861.857

$ws.Range("H61").Value = 9455.777
$ws.Range("J61").Value = 8887.75
$ws.Range("L61").Value = 8887.75
$ws.Range("N61").Value = -9583.75

$ws.Range("H62").Value = 3256.25
$ws.Range("I62").Value = 2661.5386
$ws.Range("J62").Value = 5833.3335
$ws.Range("K62").Value = 2661.5386
$ws.Range("L62").Value = 5833.3335
$ws.Range("M62").Value = -2037.5386
$ws.Range("N62").Value = -7081.3335

$ws.Range("H65").Value = 3256.25
$ws.Range("I65").Value = 2661.5386
$ws.Range("J65").Value = 5833.3335
$ws.Range("K65").Value = 13307.693
$ws.Range("L65").Value = 29166.6675
$ws.Range("M65").Value = -10187.693
$ws.Range("N65").Value = -35406.6675

$ws.Range("H68").Value = 18721.445
$ws.Range("J68").Value = 18721.445
$ws.Range("L68").Value = 18721.445
$ws.Range("N68").Value = -20219.445

$ws.Range("H71").Value = 18721.445
$ws.Range("J71").Value = 18721.445
$ws.Range("L71").Value = 56164.335
$ws.Range("N71").Value = -63652.335

$ws.Range("H99").Value = 19166.666
$ws.Range("I99").Value = 16600
$ws.Range("J99").Value = 21000
$ws.Range("K99").Value = 16600
$ws.Range("L99").Value = 21000
$ws.Range("M99").Value = -15102
$ws.Range("N99").Value = -23996

$ws.Range("H126").Value = 19166.666
$ws.Range("I126").Value = 16600
$ws.Range("J126").Value = 21000
$ws.Range("K126").Value = 49800
$ws.Range("L126").Value = 63000
$ws.Range("M126").Value = -47330
$ws.Range("N126").Value = -67940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2632.9119
$ws.Range("I133").Value = 2129.0908
$ws.Range("J133").Value = 2873.8696
$ws.Range("K133").Value = 6387.2724
$ws.Range("L133").Value = 8621.6088
$ws.Range("M133").Value = -1327.2724
$ws.Range("N133").Value = -18741.6088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 4000
$ws.Range("J21").Value = 4000
$ws.Range("L21").Value = 4000
$ws.Range("N21").Value = -4346

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws.Range("H30").Value = 4000
$ws.Range("J30").Value = 4000
$ws.Range("L30").Value = 4000
$ws.Range("N30").Value = -4210

$ws.Range("H102").Value = 6314.3335
$ws.Range("I102").Value = 8109.1816
$ws.Range("J102").Value = 1378.5
$ws.Range("K102").Value = 8109.1816
$ws.Range("L102").Value = 1378.5
$ws.Range("M102").Value = -6487.1816
$ws.Range("N102").Value = -4622.5

$ws.Range("H122").Value = 7269419.5
$ws.Range("I122").Value = 26298.85
$ws.Range("J122").Value = 55556892
$ws.Range("K122").Value = 78896.54999999999
$ws.Range("L122").Value = 166670676
$ws.Range("M122").Value = -76446.54999999999
$ws.Range("N122").Value = -166675576

$ws.Range("H126").Value = 8426.357
$ws.Range("I126").Value = 9205.75
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 27617.25
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -25147.25
$ws.Range("N126").Value = -16190

$ws.Range("H132").Value = 7203664.5
$ws.Range("I132").Value = 6040247.5
$ws.Range("J132").Value = 15153680
$ws.Range("K132").Value = 18120742.5
$ws.Range("L132").Value = 45461040
$ws.Range("M132").Value = -18118212.5
$ws.Range("N132").Value = -45466100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1469.9445
$ws.Range("I7").Value = 1366.4615
$ws.Range("J7").Value = 1739
$ws.Range("K7").Value = 1366.4615
$ws.Range("L7").Value = 1739
$ws.Range("M7").Value = -1254.4615
$ws.Range("N7").Value = -1963

$ws.Range("H40").Value = 2669.2307
$ws.Range("I40").Value = 3128.5715
$ws.Range("J40").Value = 2133.3333
$ws.Range("K40").Value = 3128.5715
$ws.Range("L40").Value = 2133.3333
$ws.Range("M40").Value = -2992.5715
$ws.Range("N40").Value = -2405.3333

$ws.Range("H122").Value = 11065990
$ws.Range("I122").Value = 1253635.1
$ws.Range("J122").Value = 66669332
$ws.Range("K122").Value = 3760905.3
$ws.Range("L122").Value = 200007996
$ws.Range("M122").Value = -3758455.3
$ws.Range("N122").Value = -200012896

$ws.Range("H126").Value = 1469.9445
$ws.Range("I126").Value = 1366.4615
$ws.Range("J126").Value = 1739
$ws.Range("K126").Value = 4099.3845
$ws.Range("L126").Value = 5217
$ws.Range("M126").Value = -1629.3845
$ws.Range("N126").Value = -10157

$ws.Range("H136").Value = 3585795
$ws.Range("I136").Value = 5051869.5
$ws.Range("J136").Value = 2056.6667
$ws.Range("K136").Value = 15155608.5
$ws.Range("L136").Value = 6170.000100000001
$ws.Range("M136").Value = -15153058.5
$ws.Range("N136").Value = -11270.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70017
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70017
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70017
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -70487

$ws.Range("H35").Value = 70017
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 70017
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 70017
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -70597

$ws.Range("H122").Value = 1693.3043
$ws.Range("I122").Value = 1412.625
$ws.Range("J122").Value = 2334.8572
$ws.Range("K122").Value = 4237.875
$ws.Range("L122").Value = 7004.571599999999
$ws.Range("M122").Value = -1787.875
$ws.Range("N122").Value = -11904.5716

$ws.Range("H136").Value = 2945.8545
$ws.Range("I136").Value = 715.2
$ws.Range("J136").Value = 8894.267
$ws.Range("K136").Value = 2145.6
$ws.Range("L136").Value = 26682.801
$ws.Range("M136").Value = 404.3999999999996
$ws.Range("N136").Value = -31782.801
